$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.916.21"
$ws.Range("E2").Value = "  +2.30%  "
$ws.Range("D3").Value = "3.035.67"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.91"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.72"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +6.91%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.029.65"
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.58"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +11.54%  "
$ws.Range("E11").Value = "  +2.90%  "
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000235"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.73"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +4.30%  "
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("D16").Value = "3.537.53"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").Value = "62.871.75"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").Value = "3.033.95"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "453.36"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.29"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.19"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.22"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +5.36%  "
$ws.Range("E26").Value = "  +3.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.39"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +3.53%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  +3.18%  "
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("E31").Value = "  +7.20%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.63"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("E35").Value = "  +4.86%  "
$ws.Range("E36").Value = "  +2.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.95"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +3.03%  "
$ws.Range("E38").Value = "  +10.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.12"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +2.43%  "
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.128"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +4.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.09"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -1.84%  "
$ws.Range("E43").Value = "  +12.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.95"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +6.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "397.55"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("E46").Value = "  +1.93%  "
$ws.Range("D47").Value = "2.728.48"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  +5.84%  "
